$d = $word.ActiveDocument

$replacements = @(
    @("148÷3=49, 1", "200÷2=100, 0"),
    @("640÷3=213, 1", "567÷8=70, 7"),
    @("281÷3=93, 2", "496÷5=99, 1"),
    @("670÷3=223, 1", "948÷3=316, 0"),
    @("435÷2=217, 1", "723÷5=144, 3"),
    @("862÷2=431, 0", "749÷3=249, 2"),
    @("167÷2=83, 1", "560÷2=280, 0"),
    @("685÷7=97, 6", "957÷3=319, 0"),
    @("944÷9=104, 8", "982÷4=245, 2"),
    @("692÷8=86, 4", "266÷6=44, 2"),
    @("606÷6=101, 0", "662÷4=165, 2"),
    @("836÷2=418, 0", "375÷5=75, 0"),
    @("442÷9=49, 1", "545÷2=272, 1"),
    @("833÷4=208, 1", "984÷8=123, 0"),
    @("257÷2=128, 1", "350÷8=43, 6"),
    @("489÷5=97, 4", "739÷4=184, 3"),
    @("737÷8=92, 1", "374÷8=46, 6"),
    @("771÷2=385, 1", "930÷2=465, 0"),
    @("433÷6=72, 1", "687÷3=229, 0"),
    @("994÷6=165, 4", "847÷2=423, 1"),
    @("838÷3=279, 1", "651÷3=217, 0"),
    @("485÷3=161, 2", "218÷9=24, 2"),
    @("225÷5=45, 0", "949÷4=237, 1"),
    @("908÷4=227, 0", "924÷3=308, 0"),
    @("165÷9=18, 3", "985÷4=246, 1")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
